$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("Add New School Department", "PASSED", "chrome"),
    @("Edit School Department", "FAILED", "chrome"),
    @("Delete The School Department", "PASSED", "chrome"),
    @("Add New School Department", "PASSED", "chrome"),
    @("Add New School Department", "PASSED", "chrome"),
    @("Edit School Department", "FAILED", "chrome"),
    @("Delete The School Department", "PASSED", "chrome"),
    @("Add New School Department", "PASSED", "chrome"),
    @("Add New School Department", "PASSED", "chrome"),
    @("Edit School Department", "FAILED", "chrome"),
    @("Add New School Department", "PASSED", "chrome"),
    @("Add New School Department", "PASSED", "chrome"),
    @("Edit School Department", "PASSED", "chrome"),
    @("Delete The School Department", "PASSED", "chrome"),
    @("Add New School Department", "PASSED", "chrome"),
    @("Edit School Department", "PASSED", "chrome"),
    @("Delete The School Department", "PASSED", "chrome"),
    @("Add the School Locations", "PASSED", "chrome"),
    @("Edit the School Locations", "PASSED", "chrome"),
    @("Add the School Locations", "PASSED", "chrome"),
    @("Edit the School Locations", "PASSED", "chrome"),
    @("Negative Test in School Locations", "PASSED", "chrome"),
    @("User Delete the School Locations", "PASSED", "chrome")
)

$startRow = 22
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $rows[$i][0]
    $ws.Cells.Item($r, 2).Value = $rows[$i][1]
    $ws.Cells.Item($r, 3).Value = $rows[$i][2]
}
